# Commit: Mon, May 11, 2020 5:14:21 PM
#
# The table on slide 6 (the "SOURCES OF FINANCE" table) gets a new
# built-in table style applied from the Table Design gallery:
#   {E588D616-27A3-4A74-A2B0-7AFBDD5B86B1}  ->  {F5B2E416-26F3-4EB1-B03C-E45A7A7B706E}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

foreach ($sh in $s.Shapes) {
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{F5B2E416-26F3-4EB1-B03C-E45A7A7B706E}")
    }
}
